$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 584.375
$ws.Range("I4").Value = 321.5
$ws.Range("K4").Value = 321.5
$ws.Range("M4").Value = -207.5
$ws.Range("H12").Value = 491.5
$ws.Range("I12").Value = 492.33334
$ws.Range("J12").Value = 489
$ws.Range("K12").Value = 492.33334
$ws.Range("L12").Value = 489
$ws.Range("M12").Value = -322.33334
$ws.Range("N12").Value = -829
$ws.Range("H17").Value = 3861.879
$ws.Range("J17").Value = 3861.879
$ws.Range("L17").Value = 11585.637
$ws.Range("N17").Value = -11921.637
$ws.Range("H18").Value = 764.7692
$ws.Range("I18").Value = 813.9167
$ws.Range("K18").Value = 813.9167
$ws.Range("M18").Value = -529.9167
$ws.Range("H80").Value = 1423.125
$ws.Range("I80").Value = 515.1111
$ws.Range("J80").Value = 1967.9333
$ws.Range("K80").Value = 1545.3333
$ws.Range("L80").Value = 5903.7999
$ws.Range("M80").Value = -547.3332999999998
$ws.Range("N80").Value = -7899.7999
$ws.Range("H83").Value = 1423.125
$ws.Range("I83").Value = 515.1111
$ws.Range("J83").Value = 1967.9333
$ws.Range("K83").Value = 4635.9999
$ws.Range("L83").Value = 17711.3997
$ws.Range("M83").Value = 356.0001000000002
$ws.Range("N83").Value = -27695.3997
$ws.Range("H92").Value = 1433.8966
$ws.Range("I92").Value = 1257.0769
$ws.Range("J92").Value = 2966.3333
$ws.Range("K92").Value = 1257.0769
$ws.Range("L92").Value = 2966.3333
$ws.Range("M92").Value = -9.076900000000023
$ws.Range("N92").Value = -5462.3333
$ws.Range("H112").Value = 1638.5333
$ws.Range("J112").Value = 1802.1538
$ws.Range("L112").Value = 5406.4614
$ws.Range("N112").Value = -7622.4614
$ws.Range("H138").Value = 2666.9092
$ws.Range("J138").Value = 3106.3928
$ws.Range("L138").Value = 9319.178400000001
$ws.Range("N138").Value = -19599.1784

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H50").Value = 2340
$ws.Range("I50").Value = 600
$ws.Range("K50").Value = 600
$ws.Range("M50").Value = 114
$ws.Range("H86").Value = 49999
$ws.Range("J86").Value = 49999
$ws.Range("L86").Value = 49999
$ws.Range("N86").Value = -52371
$ws.Range("H88").Value = 1259.579
$ws.Range("I88").Value = 1139
$ws.Range("J88").Value = 1315.2307
$ws.Range("K88").Value = 1139
$ws.Range("L88").Value = 1315.2307
$ws.Range("M88").Value = -733
$ws.Range("N88").Value = -2127.2307
$ws.Range("H89").Value = 49999
$ws.Range("J89").Value = 49999
$ws.Range("L89").Value = 149997
$ws.Range("N89").Value = -161853
$ws.Range("H91").Value = 1259.579
$ws.Range("I91").Value = 1139
$ws.Range("J91").Value = 1315.2307
$ws.Range("K91").Value = 1139
$ws.Range("L91").Value = 1315.2307
$ws.Range("M91").Value = 265
$ws.Range("N91").Value = -4123.2307
$ws.Range("H128").Value = 204725.1
$ws.Range("J128").Value = 204725.1
$ws.Range("L128").Value = 204725.1
$ws.Range("N128").Value = -214685.1

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H58").Value = 75618.336
$ws.Range("J58").Value = 75618.336
$ws.Range("L58").Value = 75618.336
$ws.Range("N58").Value = -76206.336
$ws.Range("H86").Value = 78588.08
$ws.Range("J86").Value = 334633.16
$ws.Range("L86").Value = 334633.16
$ws.Range("N86").Value = -336879.16
$ws.Range("H89").Value = 78588.08
$ws.Range("J89").Value = 334633.16
$ws.Range("L89").Value = 1673165.8
$ws.Range("N89").Value = -1684397.8
$ws.Range("H94").Value = 1512
$ws.Range("I94").Value = 1288.6666
$ws.Range("K94").Value = 1288.6666
$ws.Range("M94").Value = -837.6666
$ws.Range("H102").Value = 13618.667
$ws.Range("J102").Value = 0
$ws.Range("L102").Value = 0
$ws.Range("N102").ClearContents()
$ws.Range("H105").Value = 27790740
$ws.Range("I105").Value = 45471216
$ws.Range("J105").Value = 7132.0713
$ws.Range("K105").Value = 45471216
$ws.Range("L105").Value = 7132.0713
$ws.Range("M105").Value = -45469469
$ws.Range("N105").Value = -10626.0713
$ws.Range("H115").Value = 0
$ws.Range("J115").Value = 0
$ws.Range("L115").Value = 0
$ws.Range("N115").ClearContents()

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 1052.619
$ws.Range("I22").Value = 372.3846
$ws.Range("J22").Value = 2158
$ws.Range("K22").Value = 372.3846
$ws.Range("L22").Value = 2158
$ws.Range("M22").Value = -22.38459999999998
$ws.Range("N22").Value = -2858
$ws.Range("H23").Value = 9997.4
$ws.Range("I23").Value = 9995
$ws.Range("J23").Value = 9998
$ws.Range("K23").Value = 9995
$ws.Range("L23").Value = 9998
$ws.Range("M23").Value = -9755
$ws.Range("N23").Value = -10478
$ws.Range("H27").Value = 9997.4
$ws.Range("I27").Value = 9995
$ws.Range("J27").Value = 9998
$ws.Range("K27").Value = 9995
$ws.Range("L27").Value = 9998
$ws.Range("M27").Value = -9803
$ws.Range("N27").Value = -10382
$ws.Range("H50").Value = 47635.2
$ws.Range("J50").Value = 47635.2
$ws.Range("L50").Value = 47635.2
$ws.Range("N50").Value = -48885.2
$ws.Range("H99").Value = 3599.4375
$ws.Range("I99").Value = 3545.1333
$ws.Range("K99").Value = 3545.1333
$ws.Range("M99").Value = -2047.1333
$ws.Range("H105").Value = 1403.8
$ws.Range("J105").Value = 2011
$ws.Range("L105").Value = 2011
$ws.Range("N105").Value = -5505
$ws.Range("H109").Value = 37608.062
$ws.Range("J109").Value = 37608.062
$ws.Range("L109").Value = 37608.062
$ws.Range("N109").Value = -39688.062
$ws.Range("H117").Value = 110360.25
$ws.Range("J117").Value = 110360.25
$ws.Range("L117").Value = 110360.25
$ws.Range("N117").Value = -119538.25
$ws.Range("H126").Value = 3599.4375
$ws.Range("I126").Value = 3545.1333
$ws.Range("K126").Value = 10635.3999
$ws.Range("M126").Value = -8165.3999
$ws.Range("H132").Value = 34560.89
$ws.Range("I132").Value = 3410.9285
$ws.Range("K132").Value = 10232.7855
$ws.Range("M132").Value = -7702.7855
$ws.Range("H141").Value = 428296.28
$ws.Range("J141").Value = 466361.5
$ws.Range("L141").Value = 466361.5
$ws.Range("N141").Value = -476721.5

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 730.8
$ws.Range("I7").Value = 375
$ws.Range("J7").Value = 968
$ws.Range("K7").Value = 1125
$ws.Range("L7").Value = 2904
$ws.Range("M7").Value = -1013
$ws.Range("N7").Value = -3128
$ws.Range("H34").Value = 964.4286
$ws.Range("J34").Value = 6503
$ws.Range("L34").Value = 19509
$ws.Range("N34").Value = -19677
$ws.Range("H37").Value = 91110.625
$ws.Range("J37").Value = 91110.625
$ws.Range("L37").Value = 273331.875
$ws.Range("N37").Value = -273555.875
$ws.Range("H69").Value = 1864
$ws.Range("I69").Value = 1729.6
$ws.Range("J69").Value = 2200
$ws.Range("K69").Value = 5188.799999999999
$ws.Range("L69").Value = 6600
$ws.Range("M69").Value = -4377.799999999999
$ws.Range("N69").Value = -8222
$ws.Range("H72").Value = 1864
$ws.Range("I72").Value = 1729.6
$ws.Range("J72").Value = 2200
$ws.Range("K72").Value = 15566.4
$ws.Range("L72").Value = 19800
$ws.Range("M72").Value = -11510.4
$ws.Range("N72").Value = -27912
$ws.Range("H80").Value = 3825
$ws.Range("I80").Value = 2000
$ws.Range("J80").Value = 4085.7144
$ws.Range("K80").Value = 6000
$ws.Range("L80").Value = 12257.1432
$ws.Range("M80").Value = -5064
$ws.Range("N80").Value = -14129.1432
$ws.Range("H82").Value = 0
$ws.Range("I82").Value = 0
$ws.Range("K82").Value = 0
$ws.Range("M82").ClearContents()
$ws.Range("H83").Value = 3825
$ws.Range("I83").Value = 2000
$ws.Range("J83").Value = 4085.7144
$ws.Range("K83").Value = 18000
$ws.Range("L83").Value = 36771.4296
$ws.Range("M83").Value = -13320
$ws.Range("N83").Value = -46131.4296
$ws.Range("H85").Value = 0
$ws.Range("I85").Value = 0
$ws.Range("K85").Value = 0
$ws.Range("M85").ClearContents()
$ws.Range("H120").Value = 14497
$ws.Range("J120").Value = 14497
$ws.Range("L120").Value = 43491
$ws.Range("N120").Value = -53167
$ws.Range("H128").Value = 339724
$ws.Range("I128").Value = 339724
$ws.Range("K128").Value = 1019172
$ws.Range("M128").Value = -1014192
$ws.Range("H131").Value = 17551906
$ws.Range("J131").Value = 13314.909
$ws.Range("L131").Value = 39944.727
$ws.Range("N131").Value = -50024.727

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H15").Value = 19999.8
$ws.Range("I15").Value = 20000
$ws.Range("J15").Value = 19999.5
$ws.Range("K15").Value = 20000
$ws.Range("L15").Value = 19999.5
$ws.Range("M15").Value = -19712
$ws.Range("N15").Value = -20575.5
$ws.Range("H32").Value = 64983.332
$ws.Range("J32").Value = 82500
$ws.Range("L32").Value = 82500
$ws.Range("N32").Value = -83092
$ws.Range("H34").Value = 45000
$ws.Range("J34").Value = 45000
$ws.Range("L34").Value = 45000
$ws.Range("N34").Value = -45536
$ws.Range("H59").Value = 9999.5
$ws.Range("I59").Value = 0
$ws.Range("J59").Value = 9999.5
$ws.Range("K59").Value = 0
$ws.Range("L59").Value = 9999.5
$ws.Range("M59").ClearContents()
$ws.Range("N59").Value = -11165.5
$ws.Range("H76").Value = 45000
$ws.Range("J76").Value = 45000
$ws.Range("L76").Value = 45000
$ws.Range("N76").Value = -45630
$ws.Range("H79").Value = 45000
$ws.Range("J79").Value = 45000
$ws.Range("L79").Value = 45000
$ws.Range("N79").Value = -47184
$ws.Range("H80").Value = 3067.5
$ws.Range("I80").Value = 2279.6
$ws.Range("J80").Value = 4380.6665
$ws.Range("K80").Value = 2279.6
$ws.Range("L80").Value = 4380.6665
$ws.Range("M80").Value = -1281.6
$ws.Range("N80").Value = -6376.6665
$ws.Range("H81").Value = 19999.8
$ws.Range("I81").Value = 20000
$ws.Range("J81").Value = 19999.5
$ws.Range("K81").Value = 20000
$ws.Range("L81").Value = 19999.5
$ws.Range("M81").Value = -19002
$ws.Range("N81").Value = -21995.5
$ws.Range("H83").Value = 3067.5
$ws.Range("I83").Value = 2279.6
$ws.Range("J83").Value = 4380.6665
$ws.Range("K83").Value = 11398
$ws.Range("L83").Value = 21903.3325
$ws.Range("M83").Value = -6406
$ws.Range("N83").Value = -31887.3325
$ws.Range("H84").Value = 19999.8
$ws.Range("I84").Value = 20000
$ws.Range("J84").Value = 19999.5
$ws.Range("K84").Value = 60000
$ws.Range("L84").Value = 59998.5
$ws.Range("M84").Value = -55008
$ws.Range("N84").Value = -69982.5
$ws.Range("H118").Value = 32770
$ws.Range("J118").Value = 32770
$ws.Range("L118").Value = 32770
$ws.Range("N118").Value = -36084
$ws.Range("H129").Value = 21830
$ws.Range("J129").Value = 21830
$ws.Range("L129").Value = 21830
$ws.Range("N129").Value = -31830
$ws.Range("H132").Value = 7841.3
$ws.Range("I132").Value = 3733.1667
$ws.Range("K132").Value = 11199.5001
$ws.Range("M132").Value = -8669.500100000001

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 3649.2896
$ws.Range("J22").Value = 4477.591
$ws.Range("L22").Value = 4477.591
$ws.Range("N22").Value = -5067.591
$ws.Range("H25").Value = 16000
$ws.Range("J25").Value = 16000
$ws.Range("L25").Value = 16000
$ws.Range("N25").Value = -16460
$ws.Range("H27").Value = 3649.2896
$ws.Range("J27").Value = 4477.591
$ws.Range("L27").Value = 4477.591
$ws.Range("N27").Value = -4691.591
$ws.Range("H80").Value = 80128
$ws.Range("J80").Value = 80128
$ws.Range("L80").Value = 80128
$ws.Range("N80").Value = -82374
$ws.Range("H83").Value = 80128
$ws.Range("J83").Value = 80128
$ws.Range("L83").Value = 240384
$ws.Range("N83").Value = -251616
$ws.Range("H125").Value = 162748.5
$ws.Range("J125").Value = 162748.5
$ws.Range("L125").Value = 162748.5
$ws.Range("N125").Value = -172588.5
$ws.Range("H136").Value = 4572.838
$ws.Range("I136").Value = 4492.16
$ws.Range("J136").Value = 4740.9165
$ws.Range("K136").Value = 13476.48
$ws.Range("L136").Value = 14222.7495
$ws.Range("M136").Value = -10926.48
$ws.Range("N136").Value = -19322.7495

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 7291.4
$ws.Range("I62").Value = 13000
$ws.Range("K62").Value = 13000
$ws.Range("M62").Value = -12376
$ws.Range("H65").Value = 7291.4
$ws.Range("I65").Value = 13000
$ws.Range("K65").Value = 65000
$ws.Range("M65").Value = -61880
$ws.Range("H81").Value = 3305.1292
$ws.Range("I81").Value = 2802.1785
$ws.Range("J81").Value = 7999.3335
$ws.Range("K81").Value = 5604.357
$ws.Range("L81").Value = 15998.667
$ws.Range("M81").Value = -4543.357
$ws.Range("N81").Value = -18120.667
$ws.Range("H84").Value = 3305.1292
$ws.Range("I84").Value = 2802.1785
$ws.Range("J84").Value = 7999.3335
$ws.Range("K84").Value = 28021.785
$ws.Range("L84").Value = 79993.33499999999
$ws.Range("M84").Value = -22717.785
$ws.Range("N84").Value = -90601.33499999999
$ws.Range("H94").Value = 35338.4
$ws.Range("J94").Value = 35338.4
$ws.Range("L94").Value = 35338.4
$ws.Range("N94").Value = -37140.4
$ws.Range("H95").Value = 40184.6
$ws.Range("J95").Value = 40184.6
$ws.Range("L95").Value = 40184.6
$ws.Range("N95").Value = -45676.6
$ws.Range("H131").Value = 144664.28
$ws.Range("I131").Value = 75325
$ws.Range("J131").Value = 172400
$ws.Range("K131").Value = 75325
$ws.Range("L131").Value = 172400
$ws.Range("M131").Value = -70285
$ws.Range("N131").Value = -182480
$ws.Range("H132").Value = 5515.269
$ws.Range("I132").Value = 3926.1052
$ws.Range("J132").Value = 9828.714
$ws.Range("K132").Value = 11778.3156
$ws.Range("L132").Value = 29486.142
$ws.Range("M132").Value = -9248.3156
$ws.Range("N132").Value = -34546.142
